$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data range first so no stale values leak through
# after rows shift position.
$ws.Range("A2:H38").ClearContents()

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '85'
$ws.Range("B2").Value = 37

$ws.Range("A3").Value = 'ACPA ES'
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("A4").Value = 'ACPA HS'
$ws.Range("B4").Value = 78
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("H4").Value = 1

$ws.Range("A5").Value = 'Akron Preparatory School'
$ws.Range("B5").Value = 63
$ws.Range("C5").Value = 4
$ws.Range("E5").Value = 11
$ws.Range("G5").Value = 1

$ws.Range("A6").Value = 'Canton College Prep'
$ws.Range("B6").Value = 45
$ws.Range("C6").Value = 4

$ws.Range("A7").Value = 'Canton Harbor High School'
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("A8").Value = 'Cascade Career Prep High School'
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = 'Columbus Humanities Arts and Technol'
$ws.Range("B9").Value = 46
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("G9").Value = 2

$ws.Range("A10").Value = 'Columbus Performance Academy Shepard'
$ws.Range("B10").Value = 48
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = 5
$ws.Range("G10").Value = 4

$ws.Range("A11").Value = 'Columbus Preparatory Academy ES'
$ws.Range("C11").Value = 2

$ws.Range("A12").Value = 'Columbus Preparatory Academy HS'
$ws.Range("C12").Value = 1

$ws.Range("A13").Value = 'Columbus Preparatory and Fitness Academy ES'
$ws.Range("B13").Value = 63
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 1

$ws.Range("A14").Value = 'Columbus Preparatory and Fitness Academy HS'
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 1

$ws.Range("A15").Value = 'Dayton SMART Elementary School'
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 1

$ws.Range("A16").Value = 'Eastland Performance Academy'
$ws.Range("B16").Value = 52
$ws.Range("C16").Value = 45
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 2

$ws.Range("A17").Value = 'Fairfield Preparatory Academy'
$ws.Range("B17").Value = 35
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("G17").Value = 5

$ws.Range("A18").Value = 'Foundation Academy'
$ws.Range("B18").Value = 94
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 3

$ws.Range("A19").Value = 'Great River Connections Academy ES'
$ws.Range("B19").Value = 213
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 3

$ws.Range("A20").Value = 'Great River Connections Academy HS'
$ws.Range("B20").Value = 220
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3

$ws.Range("A21").Value = 'Lake Erie Preparatory School'
$ws.Range("B21").Value = 32
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 2

$ws.Range("A22").Value = 'Legacy Academy of Excellence'
$ws.Range("B22").Value = 17
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 1

$ws.Range("A23").Value = 'Mater Academy Preparatory'
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1

$ws.Range("A24").Value = 'Middletown Preparatory and Fitness Academy'
$ws.Range("B24").Value = 78
$ws.Range("C24").Value = 5

$ws.Range("A25").Value = 'Mt. Healthy Preparatory and Fitness'
$ws.Range("B25").Value = 50
$ws.Range("C25").Value = 3
$ws.Range("E25").Value = 2

$ws.Range("A26").Value = 'Northland Preparatory and Fitness Academy'
$ws.Range("B26").Value = 24
$ws.Range("C26").Value = 2

$ws.Range("A27").Value = 'Ohio Connections Academy ES'
$ws.Range("B27").Value = 526
$ws.Range("C27").Value = 28
$ws.Range("D27").Value = 1

$ws.Range("A28").Value = 'Ohio Connections Academy HS'
$ws.Range("B28").Value = 328
$ws.Range("C28").Value = 17
$ws.Range("D28").Value = 1

$ws.Range("A29").Value = 'Ohio Construction Academy'
$ws.Range("B29").Value = 29
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 2

$ws.Range("A30").Value = 'Ohio Digital Learning School'
$ws.Range("B30").Value = 226
$ws.Range("C30").Value = 12

$ws.Range("A31").Value = 'Ohio Virtual Academy ES'
$ws.Range("B31").Value = 1432
$ws.Range("C31").Value = 104

$ws.Range("A32").Value = 'Ohio Virtual Academy HS'
$ws.Range("B32").Value = 1431
$ws.Range("C32").Value = 57

$ws.Range("A33").Value = 'Riverside Academy'
$ws.Range("B33").Value = 35
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 6

$ws.Range("A34").Value = 'STEAM Academy of Warrensville'
$ws.Range("B34").Value = 45
$ws.Range("C34").Value = 45
$ws.Range("E34").Value = 2
$ws.Range("H34").Value = 1

$ws.Range("A35").Value = 'Sheffield Academy'
$ws.Range("B35").Value = 15
$ws.Range("C35").Value = 1
$ws.Range("E35").Value = 1
$ws.Range("G35").Value = 4

$ws.Range("A36").Value = 'Skyway Career Prep High School'
$ws.Range("B36").Value = 84
$ws.Range("C36").Value = 4
$ws.Range("E36").Value = 3

$ws.Range("A37").Value = 'South Scioto Performance Academy'
$ws.Range("B37").Value = 50
$ws.Range("C37").Value = 4
$ws.Range("E37").Value = 2

$ws.Range("A38").Value = 'South Side Academy'
$ws.Range("B38").Value = 28
$ws.Range("C38").Value = 2
$ws.Range("E38").Value = 1

$ws.Range("A39").Value = 'Springfield Preparatory and Fitness Academy'
$ws.Range("B39").Value = 53
$ws.Range("C39").Value = 45
$ws.Range("D39").Value = 1

$ws.Range("A40").Value = 'SunBridge Performance Academy'
$ws.Range("B40").Value = 55
$ws.Range("C40").Value = 4
$ws.Range("E40").Value = 2

$ws.Range("A41").Value = 'T-Squared Honors Academy North'
$ws.Range("B41").Value = 11
$ws.Range("C41").Value = 1
$ws.Range("E41").Value = 7

$ws.Range("A42").Value = 'Toledo Preparatory and Fitness Academy'
$ws.Range("B42").Value = 59
$ws.Range("C42").Value = 4
$ws.Range("E42").Value = 1

$ws.Range("A43").Value = 'Trotwood Preparatory and Fitness Academy'
$ws.Range("B43").Value = 54
$ws.Range("C43").Value = 4
$ws.Range("E43").Value = 2

$ws.Range("A44").Value = 'Unknown School'
$ws.Range("B44").Value = 21
$ws.Range("C44").Value = 1
$ws.Range("E44").Value = 1

$ws.Range("A45").Value = 'Western Toledo Preparatory'
$ws.Range("B45").Value = 4
$ws.Range("C45").Value = 4
$ws.Range("D45").Value = 2
$ws.Range("E45").Value = 3

$ws.Range("A46").Value = 'Whitehall Preparatory and Fitness Academy'
$ws.Range("B46").Value = 70
$ws.Range("C46").Value = 5
$ws.Range("E46").Value = 1

$ws.Range("A47").Value = 'Wildwood Environmental Academy ES'
$ws.Range("B47").Value = 63
$ws.Range("C47").Value = 6
$ws.Range("D47").Value = 1

$ws.Range("A48").Value = 'Wildwood Environmental Academy HS'
$ws.Range("B48").Value = 24
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = 1

Write-Host "done"